$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Movie 2"
$ws.Range("D1").Value = "Movie 3"
$ws.Range("E1").Value = "Movie 4"
$ws.Range("F1").Value = "Movie 5"

$ws.Range("A7").Value = "James"
$ws.Range("B7").Value = "223 jump street"
$ws.Range("C7").Value = "Die Hart"
$ws.Range("D7").Value = "Monkey and me"
$ws.Range("E7").Value = "I love You most"
$ws.Range("F7").Value = "Frenchie"
$ws.Range("G7").Value = "rawr"

$ws.Range("G1").Value = "Movie 6"

$ws.Range("E7").Select()
